$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 156.48
$ws.Range("I15").Value = 156.48
$ws.Range("K15").Value = 469.4399999999999
$ws.Range("M15").Value = -300.4399999999999
$ws.Range("H32").Value = 400
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""
$ws.Range("H39").Value = 201
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""
$ws.Range("H40").Value = 1345
$ws.Range("I40").Value = 1345
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1345
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1170
$ws.Range("N40").Value = ""
$ws.Range("H51").Value = 7800
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 11666.667
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 11666.667
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -12634.667
$ws.Range("H52").Value = 100500
$ws.Range("I52").Value = 1000
$ws.Range("K52").Value = 3000
$ws.Range("M52").Value = -2840
$ws.Range("H69").Value = 4804.3335
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""
$ws.Range("H70").Value = 5849297.5
$ws.Range("I70").Value = 700
$ws.Range("J70").Value = 6580372
$ws.Range("K70").Value = 2100
$ws.Range("L70").Value = 19741116
$ws.Range("M70").Value = -1830
$ws.Range("N70").Value = -19741656
$ws.Range("H72").Value = 4804.3335
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""
$ws.Range("H73").Value = 5849297.5
$ws.Range("I73").Value = 700
$ws.Range("J73").Value = 6580372
$ws.Range("K73").Value = 2100
$ws.Range("L73").Value = 19741116
$ws.Range("M73").Value = -1164
$ws.Range("N73").Value = -19742988
$ws.Range("H80").Value = 416.17648
$ws.Range("I80").Value = 414.53845
$ws.Range("J80").Value = 421.5
$ws.Range("K80").Value = 1243.61535
$ws.Range("L80").Value = 1264.5
$ws.Range("M80").Value = -245.61535
$ws.Range("N80").Value = -3260.5
$ws.Range("H83").Value = 416.17648
$ws.Range("I83").Value = 414.53845
$ws.Range("J83").Value = 421.5
$ws.Range("K83").Value = 3730.84605
$ws.Range("L83").Value = 3793.5
$ws.Range("M83").Value = 1261.15395
$ws.Range("N83").Value = -13777.5
$ws.Range("H92").Value = 448.4
$ws.Range("I92").Value = 448.4
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 448.4
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 799.6
$ws.Range("N92").Value = ""
$ws.Range("H138").Value = 3894.2083
$ws.Range("I138").Value = 1816.9333
$ws.Range("J138").Value = 4440.86
$ws.Range("K138").Value = 5450.7999
$ws.Range("L138").Value = 13322.58
$ws.Range("M138").Value = -310.7999
$ws.Range("N138").Value = -23602.58

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13175.302
$ws.Range("I32").Value = 9711.465
$ws.Range("J32").Value = 18140.133
$ws.Range("K32").Value = 9711.465
$ws.Range("L32").Value = 18140.133
$ws.Range("M32").Value = -9424.465
$ws.Range("N32").Value = -18714.133
$ws.Range("H45").Value = 1187.5
$ws.Range("I45").Value = 1150
$ws.Range("K45").Value = 1150
$ws.Range("M45").Value = -773
$ws.Range("H63").Value = 10658606
$ws.Range("I63").Value = 23087150
$ws.Range("J63").Value = 5568.5713
$ws.Range("K63").Value = 23087150
$ws.Range("L63").Value = 5568.5713
$ws.Range("M63").Value = -23086464
$ws.Range("N63").Value = -6940.5713
$ws.Range("H66").Value = 10658606
$ws.Range("I66").Value = 23087150
$ws.Range("J66").Value = 5568.5713
$ws.Range("K66").Value = 115435750
$ws.Range("L66").Value = 27842.8565
$ws.Range("M66").Value = -115432318
$ws.Range("N66").Value = -34706.85649999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 639.8
$ws.Range("I22").Value = 399.5
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 399.5
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -226.5
$ws.Range("N22").Value = -1146
$ws.Range("H134").Value = 3619.6875
$ws.Range("I134").Value = 2154.4614
$ws.Range("J134").Value = 9969
$ws.Range("K134").Value = 6463.3842
$ws.Range("L134").Value = 29907
$ws.Range("M134").Value = -3928.3842
$ws.Range("N134").Value = -34977

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2308.8071
$ws.Range("I58").Value = 1865.4615
$ws.Range("J58").Value = 6919.6
$ws.Range("K58").Value = 1865.4615
$ws.Range("L58").Value = 6919.6
$ws.Range("M58").Value = -1662.4615
$ws.Range("N58").Value = -7325.6
$ws.Range("H136").Value = 2308.8071
$ws.Range("I136").Value = 1865.4615
$ws.Range("J136").Value = 6919.6
$ws.Range("K136").Value = 5596.3845
$ws.Range("L136").Value = 20758.8
$ws.Range("M136").Value = -3046.3845
$ws.Range("N136").Value = -25858.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3968339.5
$ws.Range("I2").Value = 90
$ws.Range("K2").Value = 540
$ws.Range("M2").Value = -427
$ws.Range("H132").Value = 1612.6666
$ws.Range("I132").Value = 871.2941
$ws.Range("K132").Value = 7841.6469
$ws.Range("M132").Value = -5311.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3390.4167
$ws.Range("I132").Value = 876.4286
$ws.Range("J132").Value = 4425.5884
$ws.Range("K132").Value = 2629.2858
$ws.Range("L132").Value = 13276.7652
$ws.Range("M132").Value = -99.28579999999965
$ws.Range("N132").Value = -18336.7652

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 7203.8
$ws.Range("I9").Value = 340
$ws.Range("J9").Value = 10145.429
$ws.Range("K9").Value = 340
$ws.Range("L9").Value = 10145.429
$ws.Range("M9").Value = -116
$ws.Range("N9").Value = -10593.429
$ws.Range("H82").Value = 5794.2383
$ws.Range("I82").Value = 6607.706
$ws.Range("J82").Value = 2337
$ws.Range("K82").Value = 6607.706
$ws.Range("L82").Value = 2337
$ws.Range("M82").Value = -6246.706
$ws.Range("N82").Value = -3059
$ws.Range("H85").Value = 5794.2383
$ws.Range("I85").Value = 6607.706
$ws.Range("J85").Value = 2337
$ws.Range("K85").Value = 6607.706
$ws.Range("L85").Value = 2337
$ws.Range("M85").Value = -5359.706
$ws.Range("N85").Value = -4833
$ws.Range("H127").Value = 31683.824
$ws.Range("J127").Value = 31683.824
$ws.Range("L127").Value = 31683.824
$ws.Range("N127").Value = -41603.824

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 44995
$ws.Range("J68").Value = 44995
$ws.Range("L68").Value = 44995
$ws.Range("N68").Value = -46617
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""
$ws.Range("H71").Value = 44995
$ws.Range("J71").Value = 44995
$ws.Range("L71").Value = 134985
$ws.Range("N71").Value = -143097
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360
$ws.Range("H80").Value = 45000
$ws.Range("J80").Value = 45000
$ws.Range("L80").Value = 45000
$ws.Range("N80").Value = -46996
$ws.Range("H82").Value = 48500
$ws.Range("J82").Value = 48500
$ws.Range("L82").Value = 48500
$ws.Range("N82").Value = -49266
$ws.Range("H83").Value = 45000
$ws.Range("J83").Value = 45000
$ws.Range("L83").Value = 135000
$ws.Range("N83").Value = -144984
$ws.Range("H85").Value = 48500
$ws.Range("J85").Value = 48500
$ws.Range("L85").Value = 48500
$ws.Range("N85").Value = -51152
$ws.Range("H132").Value = 23819960
$ws.Range("I132").Value = 20881.6
$ws.Range("J132").Value = 37041668
$ws.Range("K132").Value = 62644.8
$ws.Range("L132").Value = 111125004
$ws.Range("M132").Value = -60114.8
$ws.Range("N132").Value = -111130064
